# tourney_results.xlsx — fill in the 2016 Round-of-16 (and later rounds)
# results on the "results" sheet. Setting the "Winner (enter 1 or 2)"
# column (G) for rows 54-68 is the only input needed: every other cell
# that changes (columns E/F/H/J/K/M/N/O on "results", the mirrored
# columns on "slot_results", and the rows on "results_for_sim") is
# driven by formulas that already exist in the workbook and recalculate
# automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("results")

# Game winners: 1 = the "strong seed" (column E) team won,
#               2 = the "weak seed" (column F) team won.
$winners = @{
    54 = 1   # #1 Kansas        beats #5 Maryland
    55 = 1   # #2 Villanova     beats #3 Miami FL
    56 = 1   # #1 Oregon        beats #4 Duke
    57 = 1   # #2 Oklahoma      beats #3 Texas A&M
    58 = 1   # #1 North Carolina beats #5 Indiana
    59 = 2   # #6 Notre Dame    beats #7 Wisconsin
    60 = 1   # #1 Virginia      beats #4 Iowa St
    61 = 1   # #10 Syracuse     beats #11 Gonzaga
    62 = 2   # #2 Villanova     beats #1 Kansas
    63 = 2   # #2 Oklahoma      beats #1 Oregon
    64 = 1   # #1 North Carolina beats #6 Notre Dame
    65 = 2   # #10 Syracuse     beats #1 Virginia
    66 = 1   # #1 North Carolina beats #10 Syracuse
    67 = 1   # #2 Villanova     beats #2 Oklahoma
    68 = 2   # #2 Villanova     beats #1 North Carolina (champion)
}

foreach ($row in 54..68) {
    $ws.Range("G$row").Value = $winners[$row]
}

# Reflect the final active selection recorded in the sheet (cursor sitting
# on the championship-game winner cell, with the view scrolled down so
# that row 39 is the top visible row).
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 39
$win.ScrollColumn = 1
$ws.Range("G68").Select()

# Window size the workbook was last saved at.
$win.Width = 25600
$win.Height = 16060
